$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.955.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.389.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.45%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.201"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.594"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000287"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "685.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.944.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.982.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.410.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.915"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +16.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "556.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.48%  "
$ws.Range("E34").Value = "  +0.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.620.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.12%  "
$ws.Range("E38").Value = "  +3.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0735"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.97%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.81%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0430"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.340"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.18%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.96%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.130"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.60%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.57%  "
